# Apply updated dSF (column F) values for specific rows after repulling data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    9  = -6
    12 = -7
    19 = 5
    23 = -1
    31 = -3
    33 = 3
    34 = -2
    35 = 0
    38 = 2
    39 = -3
    40 = -4
    41 = 5
    43 = 4
    44 = -1
    46 = -1
    48 = 2
    58 = -2
    61 = -3
    64 = 4
    65 = 1
    66 = 8
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
